$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 2114
$ws.Range("F5").Value = 9076
$ws.Range("F6").Value = 260
$ws.Range("F9").Value = 29
$ws.Range("F10").Value = 594
$ws.Range("F13").Value = 140
$ws.Range("F14").Value = 293
$ws.Range("F16").Value = 54
$ws.Range("F17").Value = 1493
$ws.Range("F21").Value = 1372
$ws.Range("F23").Value = 231
$ws.Range("F25").Value = 91
$ws.Range("F28").Value = 304
$ws.Range("F29").Value = 304
$ws.Range("F31").Value = 12
$ws.Range("F39").Value = 133
$ws.Range("F42").Value = 34
$ws.Range("F43").Value = 491
$ws.Range("F46").Value = 210
$ws.Range("F47").Value = 46
$ws.Range("F48").Value = 44

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 50
$ws.Range("F11").Value = 225
$ws.Range("F16").Value = 671
$ws.Range("F20").Value = 74
$ws.Range("F23").Value = 931
$ws.Range("F25").Value = 1032
$ws.Range("F26").Value = 226
$ws.Range("F29").Value = 219
$ws.Range("F31").Value = 150

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 306
$ws.Range("F7").Value = 2067
$ws.Range("F8").Value = 3110

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 9076
$ws.Range("F6").Value = 306
$ws.Range("F10").Value = 260
$ws.Range("F11").Value = 2067
$ws.Range("F12").Value = 3110
$ws.Range("F14").Value = 225
$ws.Range("F18").Value = 594
$ws.Range("F21").Value = 293
$ws.Range("F22").Value = 54
$ws.Range("F23").Value = 1493
$ws.Range("F25").Value = 1372
$ws.Range("F26").Value = 231
$ws.Range("F28").Value = 91
$ws.Range("F29").Value = 304
$ws.Range("F30").Value = 304
$ws.Range("F32").Value = 931
$ws.Range("F35").Value = 226
$ws.Range("F39").Value = 133
$ws.Range("F40").Value = 219
$ws.Range("F42").Value = 150
$ws.Range("F43").Value = 491
$ws.Range("F47").Value = 210
$ws.Range("F50").Value = 46

